$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing sentinel row's (-1) "Subs" value -- column E stays the same text.
$sentinelE = $ws.Cells.Item(69, 5).Value2

# Insert 6 blank rows where the sentinel row (69) currently sits. This shifts the sentinel
# row down to row 75, carrying its original formatting along with it.
$ws.Range("A69:E74").Insert(-4121)  # xlShiftDown

# New real submissions that now occupy rows 69-74 (what used to be occupied by the sentinel).
$newRows = @(
    @(67, "10/02/2023 18:00:33", "25 - 32", "Male",   "#6 The Boss;"),
    @(68, "10/03/2023 13:19:11", "60 - 100", "Male",   "#20 Elite Chicken & Bacon Ranch;"),
    @(69, "10/04/2023 10:32:38", "25 - 32", "Female",  "#4 Supreme Meats;"),
    @(70, "10/05/2023 13:26:56", "15 - 20", "Male",    "#6 The Boss;"),
    @(71, "10/05/2023 13:28:37", "60 - 100", "Male",   "#23 The Hotshot Italiano;"),
    @(72, "10/05/2023 13:30:06", "60 - 100", "Male",   "#6 The Boss;")
)

$r = 69
foreach ($rec in $newRows) {
    # B:E are plain text cells carrying no special number format (unlike the numeric
    # timestamps used higher up in the sheet), so clear whatever the row-insert inherited.
    $ws.Range($ws.Cells.Item($r, 2), $ws.Cells.Item($r, 5)).ClearFormats()

    $ws.Cells.Item($r, 1).Value = $rec[0]
    $ws.Cells.Item($r, 2).Value = $rec[1]
    $ws.Cells.Item($r, 3).Value = $rec[2]
    $ws.Cells.Item($r, 4).Value = $rec[3]
    $ws.Cells.Item($r, 5).Value = $rec[4]

    $ws.Cells.Item($r, 1).Font.Bold = $true
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 1).VerticalAlignment = -4160
    $ws.Cells.Item($r, 1).Borders.LineStyle = 1

    $r++
}

# The sentinel row ("currently being filled out" placeholder) is now row 75; refresh its data.
$ws.Cells.Item(75, 1).Value = -1
$ws.Cells.Item(75, 2).Value = "10/05/2023 13:30:34"
$ws.Cells.Item(75, 3).Value = "0 - 2"
$ws.Cells.Item(75, 4).Value = "Female"
$ws.Cells.Item(75, 5).Value = $sentinelE

$excel.CutCopyMode = $false
